# Update the two User-Story cells whose text changed:
#  - Row 15 (Search and Sort / "sort the list of available cupcakes"):
#    the BENEFIT cell (C15) is reworded.
#  - Row 23 (Order and checkout / "As a registered Site User"):
#    the ACTION cell (B23) is reworded (payment info -> delivery info).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C15").Value = "so that I can view the results ordered by price, date, name and category."
$ws.Range("B23").Value = "I want to be able to save my delivery information when creating an order"

# Match the author's final selection/viewport on the sheet.
$ws.Activate()
$ws.Range("C23").Select()
